$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 299.85715
$ws.Range("I2").Value = 299.85715
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 299.85715
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -186.85715
$ws.Range("N2").ClearContents()
$ws.Range("H55").Value = 494
$ws.Range("I55").Value = 585.1875
$ws.Range("J55").Value = 381.76923
$ws.Range("K55").Value = 585.1875
$ws.Range("L55").Value = 381.76923
$ws.Range("M55").Value = -371.1875
$ws.Range("N55").Value = -809.76923
$ws.Range("H62").Value = 5250
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -6748
$ws.Range("H65").Value = 5250
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -33740
$ws.Range("H70").Value = 774.9167
$ws.Range("I70").Value = 444.44446
$ws.Range("J70").Value = 1766.3334
$ws.Range("K70").Value = 1333.33338
$ws.Range("L70").Value = 5299.0002
$ws.Range("M70").Value = -1063.33338
$ws.Range("N70").Value = -5839.0002
$ws.Range("H73").Value = 774.9167
$ws.Range("I73").Value = 444.44446
$ws.Range("J73").Value = 1766.3334
$ws.Range("K73").Value = 1333.33338
$ws.Range("L73").Value = 5299.0002
$ws.Range("M73").Value = -397.33338
$ws.Range("N73").Value = -7171.0002
$ws.Range("H80").Value = 3312.88
$ws.Range("I80").Value = 960.25
$ws.Range("J80").Value = 4420
$ws.Range("K80").Value = 2880.75
$ws.Range("L80").Value = 13260
$ws.Range("M80").Value = -1882.75
$ws.Range("N80").Value = -15256
$ws.Range("H82").Value = 3052.4285
$ws.Range("I82").Value = 3052.4285
$ws.Range("K82").Value = 9157.2855
$ws.Range("M82").Value = -8751.2855
$ws.Range("H83").Value = 3312.88
$ws.Range("I83").Value = 960.25
$ws.Range("J83").Value = 4420
$ws.Range("K83").Value = 8642.25
$ws.Range("L83").Value = 39780
$ws.Range("M83").Value = -3650.25
$ws.Range("N83").Value = -49764
$ws.Range("H85").Value = 3052.4285
$ws.Range("I85").Value = 3052.4285
$ws.Range("K85").Value = 9157.2855
$ws.Range("M85").Value = -7753.2855
$ws.Range("H87").Value = 29904.762
$ws.Range("J87").Value = 29904.762
$ws.Range("L87").Value = 29904.762
$ws.Range("N87").Value = -32400.762
$ws.Range("H90").Value = 29904.762
$ws.Range("J90").Value = 29904.762
$ws.Range("L90").Value = 89714.28599999999
$ws.Range("N90").Value = -102194.286
$ws.Range("H126").Value = 49199.535
$ws.Range("J126").Value = 49199.535
$ws.Range("L126").Value = 49199.535
$ws.Range("N126").Value = -59079.535
$ws.Range("H132").Value = 14078.744
$ws.Range("I132").Value = 2046.8611
$ws.Range("K132").Value = 6140.5833
$ws.Range("M132").Value = -3610.5833
$ws.Range("H138").Value = 3314.4358
$ws.Range("I138").Value = 1797.5
$ws.Range("J138").Value = 4163.92
$ws.Range("K138").Value = 5392.5
$ws.Range("L138").Value = 12491.76
$ws.Range("M138").Value = -252.5
$ws.Range("N138").Value = -22771.76
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 21289.428
$ws.Range("I40").Value = 16250
$ws.Range("J40").Value = 23305.2
$ws.Range("K40").Value = 16250
$ws.Range("L40").Value = 23305.2
$ws.Range("M40").Value = -16074
$ws.Range("N40").Value = -23657.2
$ws.Range("H132").Value = 2555.3547
$ws.Range("I132").Value = 2634.5186
$ws.Range("J132").Value = 2021
$ws.Range("K132").Value = 7903.5558
$ws.Range("L132").Value = 6063
$ws.Range("M132").Value = -5373.5558
$ws.Range("N132").Value = -11123
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 498.83334
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 498.83334
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 498.83334
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -724.83334
$ws.Range("H33").Value = 1500
$ws.Range("I33").Value = 1500
$ws.Range("K33").Value = 1500
$ws.Range("M33").Value = -1121
$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 50000
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50368
$ws.Range("H44").Value = 9600
$ws.Range("I44").Value = 9600
$ws.Range("K44").Value = 9600
$ws.Range("M44").Value = -9158
$ws.Range("H100").Value = 78851.664
$ws.Range("J100").Value = 78851.664
$ws.Range("L100").Value = 78851.664
$ws.Range("N100").Value = -81015.664
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H132").Value = 3723
$ws.Range("I132").Value = 2140.2
$ws.Range("J132").Value = 8999
$ws.Range("K132").Value = 6420.599999999999
$ws.Range("L132").Value = 26997
$ws.Range("M132").Value = -3890.599999999999
$ws.Range("N132").Value = -32057
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1000250
$ws.Range("I51").Value = 1000250
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3000750
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -3000290
$ws.Range("N51").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5354.5
$ws.Range("I122").Value = 3959.75
$ws.Range("J122").Value = 6749.25
$ws.Range("K122").Value = 11879.25
$ws.Range("L122").Value = 20247.75
$ws.Range("M122").Value = -9429.25
$ws.Range("N122").Value = -25147.75
$ws.Range("H126").Value = 3583.1667
$ws.Range("I126").Value = 3166.3333
$ws.Range("K126").Value = 9498.999899999999
$ws.Range("M126").Value = -7028.999899999999
$ws.Range("H132").Value = 6672.75
$ws.Range("I132").Value = 5776.6216
$ws.Range("K132").Value = 17329.8648
$ws.Range("M132").Value = -14799.8648
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 36247.5
$ws.Range("I45").Value = 35000
$ws.Range("J45").Value = 37495
$ws.Range("K45").Value = 35000
$ws.Range("L45").Value = 37495
$ws.Range("M45").Value = -34593
$ws.Range("N45").Value = -38309
$ws.Range("H47").Value = 34493.5
$ws.Range("J47").Value = 34495
$ws.Range("L47").Value = 34495
$ws.Range("N47").Value = -35475
$ws.Range("H48").Value = 37495
$ws.Range("J48").Value = 37495
$ws.Range("L48").Value = 37495
$ws.Range("N48").Value = -38817
$ws.Range("H52").Value = 34493.5
$ws.Range("J52").Value = 34495
$ws.Range("L52").Value = 34495
$ws.Range("N52").Value = -34961
$ws.Range("H122").Value = 2546.65
$ws.Range("I122").Value = 2563
$ws.Range("K122").Value = 7689
$ws.Range("M122").Value = -5239
$ws.Range("H137").Value = 39573.47
$ws.Range("J137").Value = 39613.25
$ws.Range("L137").Value = 39613.25
$ws.Range("N137").Value = -49813.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18849.5
$ws.Range("J41").Value = 18849.5
$ws.Range("L41").Value = 18849.5
$ws.Range("N41").Value = -19629.5
$ws.Range("H49").Value = 25602
$ws.Range("I49").Value = 21505
$ws.Range("K49").Value = 21505
$ws.Range("M49").Value = -21275
$ws.Range("H107").Value = 491.33334
$ws.Range("I107").Value = 491.33334
$ws.Range("K107").Value = 1474.00002
$ws.Range("M107").Value = 445.9999800000001
